$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 6999.2  # H33 (was 7858.231)
$ws.Cells.Item(33, 9).Value = 79.833336  # I33 (was 93.85714)
$ws.Cells.Item(33, 10).Value = 11612.111  # J33 (was 16916.666)
$ws.Cells.Item(33, 11).Value = 79.833336  # K33 (was 93.85714)
$ws.Cells.Item(33, 12).Value = 11612.111  # L33 (was 16916.666)
$ws.Cells.Item(33, 13).Value = 149.166664  # M33 (was 135.14286)
$ws.Cells.Item(33, 14).Value = -12070.111  # N33 (was -17374.666)
$ws.Cells.Item(116, 8).Value = 7786  # H116 (was 7990.909)
$ws.Cells.Item(116, 9).Value = 8876.471  # I116 (was 9337.5)
$ws.Cells.Item(116, 10).Value = 4696.3335  # J116 (was 4400)
$ws.Cells.Item(116, 11).Value = 8876.471  # K116 (was 9337.5)
$ws.Cells.Item(116, 12).Value = 4696.3335  # L116 (was 4400)
$ws.Cells.Item(116, 13).Value = -5434.471  # M116 (was -5895.5)
$ws.Cells.Item(116, 14).Value = -11580.3335  # N116 (was -11284)
$ws.Cells.Item(133, 8).Value = 49975  # H133 (was 50000)
$ws.Cells.Item(133, 10).Value = 49975  # J133 (was 50000)
$ws.Cells.Item(133, 12).Value = 49975  # L133 (was 50000)
$ws.Cells.Item(133, 14).Value = -60095  # N133 (was -60120)
$ws.Cells.Item(134, 8).Value = 45000  # H134 (was 44000)
$ws.Cells.Item(134, 10).Value = 45000  # J134 (was 44000)
$ws.Cells.Item(134, 12).Value = 45000  # L134 (was 44000)
$ws.Cells.Item(134, 14).Value = -55140  # N134 (was -54140)
$ws.Cells.Item(136, 8).Value = 42000  # H136 (was 0)
$ws.Cells.Item(136, 10).Value = 42000  # J136 (was 0)
$ws.Cells.Item(136, 12).Value = 42000  # L136 (was 0)
$ws.Cells.Item(136, 14).Value = -52200  # N136 (new cell)
$ws.Cells.Item(137, 8).Value = 962.8919  # H137 (was 1025.9697)
$ws.Cells.Item(137, 9).Value = 835.6111  # I137 (was 878.8125)
$ws.Cells.Item(137, 10).Value = 1083.4736  # J137 (was 1164.4706)
$ws.Cells.Item(137, 11).Value = 2506.8333  # K137 (was 2636.4375)
$ws.Cells.Item(137, 12).Value = 3250.4208  # L137 (was 3493.4118)
$ws.Cells.Item(137, 13).Value = 43.16670000000022  # M137 (was -86.4375)
$ws.Cells.Item(137, 14).Value = -8350.4208  # N137 (was -8593.4118)
$ws.Cells.Item(139, 8).Value = 70290  # H139 (was 70393.336)
$ws.Cells.Item(139, 10).Value = 70290  # J139 (was 70393.336)
$ws.Cells.Item(139, 12).Value = 70290  # L139 (was 70393.336)
$ws.Cells.Item(139, 14).Value = -80570  # N139 (was -80673.336)

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3409.6553  # H61 (was 3287.7188)
$ws.Cells.Item(61, 9).Value = 3602.7917  # I61 (was 3479.16)
$ws.Cells.Item(61, 10).Value = 2482.6  # J61 (was 2604)
$ws.Cells.Item(61, 11).Value = 3602.7917  # K61 (was 3479.16)
$ws.Cells.Item(61, 12).Value = 2482.6  # L61 (was 2604)
$ws.Cells.Item(61, 13).Value = -3390.7917  # M61 (was -3267.16)
$ws.Cells.Item(61, 14).Value = -2906.6  # N61 (was -3028)
$ws.Cells.Item(74, 8).Value = 2682.1155  # H74 (was 2853.6)
$ws.Cells.Item(74, 9).Value = 2506.2666  # I74 (was 2813.2307)
$ws.Cells.Item(74, 10).Value = 2921.9092  # J74 (was 2928.5715)
$ws.Cells.Item(74, 11).Value = 2506.2666  # K74 (was 2813.2307)
$ws.Cells.Item(74, 12).Value = 2921.9092  # L74 (was 2928.5715)
$ws.Cells.Item(74, 13).Value = -1632.2666  # M74 (was -1939.2307)
$ws.Cells.Item(74, 14).Value = -4669.9092  # N74 (was -4676.5715)
$ws.Cells.Item(77, 8).Value = 2682.1155  # H77 (was 2853.6)
$ws.Cells.Item(77, 9).Value = 2506.2666  # I77 (was 2813.2307)
$ws.Cells.Item(77, 10).Value = 2921.9092  # J77 (was 2928.5715)
$ws.Cells.Item(77, 11).Value = 12531.333  # K77 (was 14066.1535)
$ws.Cells.Item(77, 12).Value = 14609.546  # L77 (was 14642.8575)
$ws.Cells.Item(77, 13).Value = -8163.332999999999  # M77 (was -9698.1535)
$ws.Cells.Item(77, 14).Value = -23345.546  # N77 (was -23378.8575)
$ws.Cells.Item(132, 8).Value = 1809.5319  # H132 (was 1790)
$ws.Cells.Item(132, 9).Value = 1297.8611  # I132 (was 1448.037)
$ws.Cells.Item(132, 10).Value = 3484.0908  # J132 (was 2209.682)
$ws.Cells.Item(132, 11).Value = 3893.5833  # K132 (was 4344.111)
$ws.Cells.Item(132, 12).Value = 10452.2724  # L132 (was 6629.045999999999)
$ws.Cells.Item(132, 13).Value = -1363.5833  # M132 (was -1814.111)
$ws.Cells.Item(132, 14).Value = -15512.2724  # N132 (was -11689.046)
$ws.Cells.Item(133, 8).Value = 72526.375  # H133 (was 75751.57000000001)
$ws.Cells.Item(133, 10).Value = 72526.375  # J133 (was 75751.57000000001)
$ws.Cells.Item(133, 12).Value = 72526.375  # L133 (was 75751.57000000001)
$ws.Cells.Item(133, 14).Value = -77586.375  # N133 (was -80811.57000000001)
$ws.Cells.Item(135, 8).Value = 44976.332  # H135 (was 69429)
$ws.Cells.Item(135, 10).Value = 44976.332  # J135 (was 69429)
$ws.Cells.Item(135, 12).Value = 44976.332  # L135 (was 69429)
$ws.Cells.Item(135, 14).Value = -55116.332  # N135 (was -79569)
$ws.Cells.Item(136, 8).Value = 3409.6553  # H136 (was 3287.7188)
$ws.Cells.Item(136, 9).Value = 3602.7917  # I136 (was 3479.16)
$ws.Cells.Item(136, 10).Value = 2482.6  # J136 (was 2604)
$ws.Cells.Item(136, 11).Value = 10808.3751  # K136 (was 10437.48)
$ws.Cells.Item(136, 12).Value = 7447.799999999999  # L136 (was 7812)
$ws.Cells.Item(136, 13).Value = -8258.375100000001  # M136 (was -7887.48)
$ws.Cells.Item(136, 14).Value = -12547.8  # N136 (was -12912)
$ws.Cells.Item(138, 8).Value = 60825.715  # H138 (was 61196.668)
$ws.Cells.Item(138, 10).Value = 60825.715  # J138 (was 61196.668)
$ws.Cells.Item(138, 12).Value = 60825.715  # L138 (was 61196.668)
$ws.Cells.Item(138, 14).Value = -71105.715  # N138 (was -71476.66800000001)
$ws.Cells.Item(139, 8).Value = 54744.285  # H139 (was 56172.855)
$ws.Cells.Item(139, 10).Value = 54744.285  # J139 (was 56172.855)
$ws.Cells.Item(139, 12).Value = 54744.285  # L139 (was 56172.855)
$ws.Cells.Item(139, 14).Value = -65024.285  # N139 (was -66452.85500000001)
$ws.Cells.Item(140, 8).Value = 88538.164  # H140 (was 89288.164)
$ws.Cells.Item(140, 10).Value = 88538.164  # J140 (was 89288.164)
$ws.Cells.Item(140, 12).Value = 88538.164  # L140 (was 89288.164)
$ws.Cells.Item(140, 14).Value = -98898.164  # N140 (was -99648.164)
$ws.Cells.Item(141, 8).Value = 60235.266  # H141 (was 56632.9)
$ws.Cells.Item(141, 10).Value = 60235.266  # J141 (was 56632.9)
$ws.Cells.Item(141, 12).Value = 60235.266  # L141 (was 56632.9)
$ws.Cells.Item(141, 14).Value = -70595.266  # N141 (was -66992.89999999999)

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 399.32144  # H80 (was 642.2)
$ws.Cells.Item(80, 9).Value = 612.9  # I80 (was 1169.2)
$ws.Cells.Item(80, 10).Value = 280.66666  # J80 (was 466.53333)
$ws.Cells.Item(80, 11).Value = 612.9  # K80 (was 1169.2)
$ws.Cells.Item(80, 12).Value = 280.66666  # L80 (was 466.53333)
$ws.Cells.Item(80, 13).Value = 385.1  # M80 (was -171.2)
$ws.Cells.Item(80, 14).Value = -2276.66666  # N80 (was -2462.53333)
$ws.Cells.Item(81, 8).Value = 47350  # H81 (was 48542.855)
$ws.Cells.Item(81, 10).Value = 47350  # J81 (was 48542.855)
$ws.Cells.Item(81, 12).Value = 47350  # L81 (was 48542.855)
$ws.Cells.Item(81, 14).Value = -49472  # N81 (was -50664.855)
$ws.Cells.Item(83, 8).Value = 399.32144  # H83 (was 642.2)
$ws.Cells.Item(83, 9).Value = 612.9  # I83 (was 1169.2)
$ws.Cells.Item(83, 10).Value = 280.66666  # J83 (was 466.53333)
$ws.Cells.Item(83, 11).Value = 3064.5  # K83 (was 5846)
$ws.Cells.Item(83, 12).Value = 1403.3333  # L83 (was 2332.66665)
$ws.Cells.Item(83, 13).Value = 1927.5  # M83 (was -854)
$ws.Cells.Item(83, 14).Value = -11387.3333  # N83 (was -12316.66665)
$ws.Cells.Item(84, 8).Value = 47350  # H84 (was 48542.855)
$ws.Cells.Item(84, 10).Value = 47350  # J84 (was 48542.855)
$ws.Cells.Item(84, 12).Value = 142050  # L84 (was 145628.565)
$ws.Cells.Item(84, 14).Value = -152658  # N84 (was -156236.565)
$ws.Cells.Item(132, 8).Value = 50755  # H132 (was 50780)
$ws.Cells.Item(132, 10).Value = 50755  # J132 (was 50780)
$ws.Cells.Item(132, 12).Value = 50755  # L132 (was 50780)
$ws.Cells.Item(132, 14).Value = -60875  # N132 (was -60900)
$ws.Cells.Item(134, 8).Value = 3466.3433  # H134 (was 3834.9)
$ws.Cells.Item(134, 9).Value = 975.6667  # I134 (was 1081.6904)
$ws.Cells.Item(134, 10).Value = 9758.579  # J134 (was 10259.056)
$ws.Cells.Item(134, 11).Value = 2927.0001  # K134 (was 3245.0712)
$ws.Cells.Item(134, 12).Value = 29275.737  # L134 (was 30777.168)
$ws.Cells.Item(134, 13).Value = -392.0001000000002  # M134 (was -710.0711999999999)
$ws.Cells.Item(134, 14).Value = -34345.737  # N134 (was -35847.16800000001)
$ws.Cells.Item(135, 8).Value = 41666.668  # H135 (was 45000)
$ws.Cells.Item(135, 10).Value = 41666.668  # J135 (was 45000)
$ws.Cells.Item(135, 12).Value = 41666.668  # L135 (was 45000)
$ws.Cells.Item(135, 14).Value = -51806.668  # N135 (was -55140)
$ws.Cells.Item(137, 8).Value = 51344.45  # H137 (was 70958.97)
$ws.Cells.Item(137, 10).Value = 51344.45  # J137 (was 70958.97)
$ws.Cells.Item(137, 12).Value = 51344.45  # L137 (was 70958.97)
$ws.Cells.Item(137, 14).Value = -61544.45  # N137 (was -81158.97)
$ws.Cells.Item(140, 8).Value = 74866.664  # H140 (was 87325)
$ws.Cells.Item(140, 10).Value = 74866.664  # J140 (was 87325)
$ws.Cells.Item(140, 12).Value = 74866.664  # L140 (was 87325)
$ws.Cells.Item(140, 14).Value = -85226.664  # N140 (was -97685)

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(25, 8).Value = 3000  # H25 (was 3370.6667)
$ws.Cells.Item(25, 9).Value = 0  # I25 (was 200)
$ws.Cells.Item(25, 10).Value = 3000  # J25 (was 4004.8)
$ws.Cells.Item(25, 11).Value = 0  # K25 (was 200)
$ws.Cells.Item(25, 12).Value = 3000  # L25 (was 4004.8)
$ws.Cells.Item(25, 13).ClearContents()  # M25 (was -26)
$ws.Cells.Item(25, 14).Value = -3348  # N25 (was -4352.8)
$ws.Cells.Item(31, 8).Value = 4120.456  # H31 (was 4536.625)
$ws.Cells.Item(31, 9).Value = 3582.0344  # I31 (was 3700.3215)
$ws.Cells.Item(31, 10).Value = 4520.8203  # J31 (was 5068.8184)
$ws.Cells.Item(31, 11).Value = 3582.0344  # K31 (was 3700.3215)
$ws.Cells.Item(31, 12).Value = 4520.8203  # L31 (was 5068.8184)
$ws.Cells.Item(31, 13).Value = -3287.0344  # M31 (was -3405.3215)
$ws.Cells.Item(31, 14).Value = -5110.8203  # N31 (was -5658.8184)
$ws.Cells.Item(34, 8).Value = 4120.456  # H34 (was 4536.625)
$ws.Cells.Item(34, 9).Value = 3582.0344  # I34 (was 3700.3215)
$ws.Cells.Item(34, 10).Value = 4520.8203  # J34 (was 5068.8184)
$ws.Cells.Item(34, 11).Value = 3582.0344  # K34 (was 3700.3215)
$ws.Cells.Item(34, 12).Value = 4520.8203  # L34 (was 5068.8184)
$ws.Cells.Item(34, 13).Value = -3380.0344  # M34 (was -3498.3215)
$ws.Cells.Item(34, 14).Value = -4924.8203  # N34 (was -5472.8184)
$ws.Cells.Item(134, 8).Value = 1022.7879  # H134 (was 1288.381)
$ws.Cells.Item(134, 9).Value = 832.8333  # I134 (was 1033.2858)
$ws.Cells.Item(134, 10).Value = 1529.3334  # J134 (was 1798.5714)
$ws.Cells.Item(134, 11).Value = 2498.4999  # K134 (was 3099.8574)
$ws.Cells.Item(134, 12).Value = 4588.0002  # L134 (was 5395.7142)
$ws.Cells.Item(134, 13).Value = 36.5001000000002  # M134 (was -564.8574000000003)
$ws.Cells.Item(134, 14).Value = -9658.0002  # N134 (was -10465.7142)
$ws.Cells.Item(135, 8).Value = 43830.77  # H135 (was 45416.668)
$ws.Cells.Item(135, 10).Value = 43830.77  # J135 (was 45416.668)
$ws.Cells.Item(135, 12).Value = 43830.77  # L135 (was 45416.668)
$ws.Cells.Item(135, 14).Value = -53970.77  # N135 (was -55556.668)
$ws.Cells.Item(138, 8).Value = 47500  # H138 (was 50000)
$ws.Cells.Item(138, 10).Value = 47500  # J138 (was 50000)
$ws.Cells.Item(138, 12).Value = 47500  # L138 (was 50000)
$ws.Cells.Item(138, 14).Value = -57780  # N138 (was -60280)

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 3752.5  # H18 (was 4505)
$ws.Cells.Item(18, 9).Value = 3752.5  # I18 (was 4505)
$ws.Cells.Item(18, 11).Value = 3752.5  # K18 (was 4505)
$ws.Cells.Item(18, 13).Value = -3459.5  # M18 (was -4212)
$ws.Cells.Item(133, 8).Value = 34466  # H133 (was 35645)
$ws.Cells.Item(133, 10).Value = 34466  # J133 (was 35645)
$ws.Cells.Item(133, 12).Value = 34466  # L133 (was 35645)
$ws.Cells.Item(133, 14).Value = -44586  # N133 (was -45765)
$ws.Cells.Item(135, 8).Value = 47800  # H135 (was 0)
$ws.Cells.Item(135, 10).Value = 47800  # J135 (was 0)
$ws.Cells.Item(135, 12).Value = 47800  # L135 (was 0)
$ws.Cells.Item(135, 14).Value = -57940  # N135 (new cell)
$ws.Cells.Item(138, 8).Value = 61708.332  # H138 (was 64875)
$ws.Cells.Item(138, 10).Value = 61708.332  # J138 (was 64875)
$ws.Cells.Item(138, 12).Value = 61708.332  # L138 (was 64875)
$ws.Cells.Item(138, 14).Value = -71988.33199999999  # N138 (was -75155)
$ws.Cells.Item(139, 8).Value = 65750  # H139 (was 66500)
$ws.Cells.Item(139, 10).Value = 65750  # J139 (was 66500)
$ws.Cells.Item(139, 12).Value = 65750  # L139 (was 66500)
$ws.Cells.Item(139, 14).Value = -76030  # N139 (was -76780)
$ws.Cells.Item(140, 8).Value = 99744.5  # H140 (was 99754.5)
$ws.Cells.Item(140, 10).Value = 99744.5  # J140 (was 99754.5)
$ws.Cells.Item(140, 12).Value = 99744.5  # L140 (was 99754.5)
$ws.Cells.Item(140, 14).Value = -110104.5  # N140 (was -110114.5)

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2614.0962  # H136 (was 2583.5283)
$ws.Cells.Item(136, 9).Value = 1704.4359  # I136 (was 1686.675)
$ws.Cells.Item(136, 11).Value = 5113.307699999999  # K136 (was 5060.025)
$ws.Cells.Item(136, 13).Value = -2563.307699999999  # M136 (was -2510.025)
$ws.Cells.Item(138, 8).Value = 65259.5  # H138 (was 64679.57)
$ws.Cells.Item(138, 10).Value = 65259.5  # J138 (was 64679.57)
$ws.Cells.Item(138, 12).Value = 65259.5  # L138 (was 64679.57)
$ws.Cells.Item(138, 14).Value = -75539.5  # N138 (was -74959.57000000001)
$ws.Cells.Item(139, 8).Value = 69900  # H139 (was 79800)
$ws.Cells.Item(139, 10).Value = 69900  # J139 (was 79800)
$ws.Cells.Item(139, 12).Value = 69900  # L139 (was 79800)
$ws.Cells.Item(139, 14).Value = -80180  # N139 (was -90080)
$ws.Cells.Item(140, 8).Value = 78476.336  # H140 (was 58590)
$ws.Cells.Item(140, 10).Value = 78476.336  # J140 (was 58590)
$ws.Cells.Item(140, 12).Value = 78476.336  # L140 (was 58590)
$ws.Cells.Item(140, 14).Value = -88836.336  # N140 (was -68950)
$ws.Cells.Item(141, 8).Value = 67033.336  # H141 (was 65785)
$ws.Cells.Item(141, 10).Value = 67033.336  # J141 (was 65785)
$ws.Cells.Item(141, 12).Value = 67033.336  # L141 (was 65785)
$ws.Cells.Item(141, 14).Value = -77393.336  # N141 (was -76145)

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(138, 8).Value = 49525  # H138 (was 49775)
$ws.Cells.Item(138, 10).Value = 49525  # J138 (was 49775)
$ws.Cells.Item(138, 12).Value = 49525  # L138 (was 49775)
$ws.Cells.Item(138, 14).Value = -59805  # N138 (was -60055)
$ws.Cells.Item(139, 8).Value = 53942.855  # H139 (was 54683.332)
$ws.Cells.Item(139, 10).Value = 53942.855  # J139 (was 54683.332)
$ws.Cells.Item(139, 12).Value = 53942.855  # L139 (was 54683.332)
$ws.Cells.Item(139, 14).Value = -64222.855  # N139 (was -64963.332)
$ws.Cells.Item(141, 8).Value = 60128.57  # H141 (was 63833.332)
$ws.Cells.Item(141, 10).Value = 62150  # J141 (was 67000)
$ws.Cells.Item(141, 12).Value = 62150  # L141 (was 67000)
$ws.Cells.Item(141, 14).Value = -72510  # N141 (was -77360)
